$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force column D to text format so numeric-looking price strings
# (e.g. "129.44") are preserved as text, matching the source data style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.934.39"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "3.412.58"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "409.74"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").Value = "129.44"
$ws.Range("E6").Value = "  +0.57%  "

$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "0.724"
$ws.Range("E9").Value = "  -1.98%  "

$ws.Range("E10").Value = "  -4.69%  "

$ws.Range("D11").Value = "42.97"
$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("D12").Value = "9.11"
$ws.Range("E12").Value = "  +2.34%  "

$ws.Range("D13").Value = "3.956.14"
$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("E14").Value = "  +0.09%  "

$ws.Range("D15").Value = "0.0000209"
$ws.Range("E15").Value = "  -2.83%  "

$ws.Range("D16").Value = "20.86"
$ws.Range("E16").Value = "  -2.00%  "

$ws.Range("D17").Value = "3.409.84"
$ws.Range("E17").Value = "  -0.70%  "

$ws.Range("D18").Value = "12.52"
$ws.Range("E18").Value = "  +0.95%  "

$ws.Range("D19").Value = "1.09"
$ws.Range("E19").Value = "  +1.53%  "

$ws.Range("D20").Value = "61.808.89"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").Value = "478.09"
$ws.Range("E21").Value = "  +19.15%  "

$ws.Range("D22").Value = "91.27"
$ws.Range("E22").Value = "  +1.45%  "

$ws.Range("D23").Value = "3.27"
$ws.Range("E23").Value = "  +2.59%  "

$ws.Range("D24").Value = "13.30"
$ws.Range("E24").Value = "  -0.73%  "

$ws.Range("D25").Value = "3.33"
$ws.Range("E25").Value = "  +3.64%  "

$ws.Range("D26").Value = "33.98"
$ws.Range("E26").Value = "  +2.99%  "

$ws.Range("D27").Value = "9.06"
$ws.Range("E27").Value = "  +4.41%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").Value = "7.79"
$ws.Range("E29").Value = "  +2.39%  "

$ws.Range("D30").Value = "2.76"
$ws.Range("E30").Value = "  +1.30%  "

$ws.Range("D31").Value = "11.94"
$ws.Range("E31").Value = "  +0.38%  "

$ws.Range("D32").Value = "0.166"
$ws.Range("E32").Value = "  -3.36%  "

$ws.Range("E33").Value = "  -4.75%  "

$ws.Range("D34").Value = "41.58"
$ws.Range("E34").Value = "  -4.74%  "

$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").Value = "57.30"
$ws.Range("E36").Value = "  +6.33%  "

$ws.Range("D37").Value = "0.0489"
$ws.Range("E37").Value = "  -3.24%  "

$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("D39").Value = "150.39"
$ws.Range("E39").Value = "  +5.89%  "

$ws.Range("D40").Value = "3.40"
$ws.Range("E40").Value = "  +0.13%  "

$ws.Range("D41").Value = "0.135"
$ws.Range("E41").Value = "  +2.40%  "

$ws.Range("D42").Value = "0.320"
$ws.Range("E42").Value = "  +2.06%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "2.93"
$ws.Range("E43").Value = "  +0.66%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "2.10"
$ws.Range("E44").Value = "  +5.98%  "

$ws.Range("D45").Value = "2.62"
$ws.Range("E45").Value = "  +8.66%  "

$ws.Range("D46").Value = "4.19"
$ws.Range("E46").Value = "  +3.23%  "

$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").Value = "16.48"
$ws.Range("E47").Value = "  -1.23%  "

$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  +19.66%  "

$ws.Range("E49").Value = "  +15.58%  "

$ws.Range("D50").Value = "22.49"
$ws.Range("E50").Value = "  +2.58%  "

$ws.Range("D51").Value = "116.41"
$ws.Range("E51").Value = "  +20.28%  "

# Restore default (General/Normal) style so only the cell content differs,
# not the formatting -- matches original workbook styling.
$ws.Range("D2:D51").Style = "Normal"
